$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "_old" / "_new" column-header suffixes to the respective
# input-file format versions ("_FV2404" / "_FV2410").
for ($col = 1; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2404")
}
for ($col = 12; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2410")
}

# Turn the header row + data range into an Excel Table (ListObject) so the
# new header names are also exposed as the table's column names.
$range = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row (pane split after row 1).
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
